$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new blog post (ser: 149 - Surah Raad Verse 24-25) was added to the content
# calendar. The existing "blog" slot strings in the last row (row 11) cascade
# up by one serial number: the old "ser: 146" slot now shows "ser: 147", the
# old "ser: 147" slot now shows "ser: 148", and the old "ser: 148" slot now
# shows the brand new "ser: 149".
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 147"
$ws.Range("D11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 148"
$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 149"

# Reflect the author's last click/selection on the refreshed cell.
$ws.Range("I11").Select()
